# 2021 Buccaneers Team Data -- log the divisional-round game and refresh the
# running, simulated-season totals (per commit message: "Logged 2021
# divisional round, simulated season from conference round").
#
# The per-game logs on sheet YDS (rush/pass yards for/against) and sheet ST
# (kickoff/punt/FG distance logs) are simple space-separated running lists --
# one new number (or handful of numbers) appended per game. The season
# aggregate sheets (OFF/DEF/ST/TURNS/PEN) get their running totals bumped by
# the new game's counts.

$wb = $excel.ActiveWorkbook

function Append-Numbers {
    param($SheetName, $CellRef, $NewNumbers)
    $ws = $wb.Worksheets.Item($SheetName)
    $cell = $ws.Range($CellRef)
    $existing = [string]$cell.Text
    $addition = ($NewNumbers -join " ")
    if ($existing -and $existing.Length -gt 0) {
        $cell.Value = "$existing $addition"
    } else {
        $cell.Value = $addition
    }
}

# ---------------------------------------------------------------------
# YDS sheet: append this game's per-drive rush/pass yards for (B) and
# against (C) to the running season log.
# ---------------------------------------------------------------------
Append-Numbers "YDS" "B2" @("11","10","2","1","1","3","7","0","-1","1","4","3")
Append-Numbers "YDS" "B3" @("13","29","5","3","9","4","21","8","19","1","5","7","-2","42","11","1","3","9","13","5","9","6","4","-3","19","3","55","5","9")
Append-Numbers "YDS" "C2" @("1","4","9","2","0","-3","1","2","3","7","0","5","2","-3","2","0","1","0","3","15","3","8","-3","1","4","4","-1","3","3")
Append-Numbers "YDS" "C3" @("6","9","17","20","29","9","9","7","70","6","5","6","4","7","-4","11","7","11","17","13","11","4","10","3","18","0","20","44")

# ---------------------------------------------------------------------
# OFF sheet: season offensive totals after the divisional-round game
# (row 2 = Home, row 3 = Road).
# ---------------------------------------------------------------------
$off = $wb.Worksheets.Item("OFF")
$off.Range("C2").Value = 224
$off.Range("D2").Value = 15
$off.Range("F2").Value = 31
$off.Range("G2").Value = 60
$off.Range("J2").Value = 34
$off.Range("L2").Value = 387
$off.Range("M2").Value = 258
$off.Range("N2").Value = 15
$off.Range("O2").Value = 21
$off.Range("P2").Value = 11
$off.Range("Q2").Value = 596

$off.Range("C3").Value = 260
$off.Range("E3").Value = 32
$off.Range("F3").Value = 185
$off.Range("G3").Value = 77
$off.Range("H3").Value = 20
$off.Range("I3").Value = 88
$off.Range("J3").Value = 81
$off.Range("N3").Value = 16

# ---------------------------------------------------------------------
# DEF sheet: season defensive totals after the divisional-round game.
# ---------------------------------------------------------------------
$def = $wb.Worksheets.Item("DEF")
$def.Range("C2").Value = 204
$def.Range("D2").Value = 8
$def.Range("E2").Value = 11
$def.Range("F2").Value = 57
$def.Range("G2").Value = 54
$def.Range("H2").Value = 8
$def.Range("J2").Value = 26
$def.Range("L2").Value = 386
$def.Range("M2").Value = 261
$def.Range("Q2").Value = 598

$def.Range("B3").Value = 8
$def.Range("C3").Value = 268
$def.Range("E3").Value = 46
$def.Range("F3").Value = 150
$def.Range("G3").Value = 53
$def.Range("H3").Value = 32
$def.Range("I3").Value = 80
$def.Range("J3").Value = 81
$def.Range("N3").Value = 20

# ---------------------------------------------------------------------
# ST sheet: season special-teams totals (row 2/3) plus the per-kick /
# per-punt distance logs (row 3-6) get this game's kicks appended.
# ---------------------------------------------------------------------
$st = $wb.Worksheets.Item("ST")
$st.Range("B2").Value = 115
$st.Range("D2").Value = 69
$st.Range("F2").Value = 471
$st.Range("G2").Value = 453
$st.Range("J2").Value = 190
$st.Range("K2").Value = 180
$st.Range("L2").Value = 121
$st.Range("M2").Value = 91
$st.Range("B3").Value = 83

Append-Numbers "ST" "D3" @("57","51","55","49")
Append-Numbers "ST" "B4" @("64","63","64")
Append-Numbers "ST" "D4" @("0","6","21","33")
Append-Numbers "ST" "B5" @("0","0","19")
Append-Numbers "ST" "D5" @("0","0","5","2")
Append-Numbers "ST" "B6" @("20","24","4","14")

# ---------------------------------------------------------------------
# TURNS sheet: season turnover totals.
# ---------------------------------------------------------------------
$turns = $wb.Worksheets.Item("TURNS")
$turns.Range("B2").Value = 4
$turns.Range("D2").Value = 7
$turns.Range("E2").Value = 13
$turns.Range("E3").Value = 10

# ---------------------------------------------------------------------
# PEN sheet: season penalty totals.
# ---------------------------------------------------------------------
$pen = $wb.Worksheets.Item("PEN")
$pen.Range("D4").Value = 15

Write-Host "2021 divisional round logged."
